$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DatosGenerales")
$ws.Range("B6").Value = "24"
